$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1012.375
$ws.Range("I2").Value = 2533.3333
$ws.Range("J2").Value = 99.8
$ws.Range("K2").Value = 2533.3333
$ws.Range("L2").Value = 99.8
$ws.Range("M2").Value = -2420.3333
$ws.Range("N2").Value = -325.8
$ws.Range("H33").Value = 999.5
$ws.Range("I33").Value = 999.5
$ws.Range("K33").Value = 999.5
$ws.Range("M33").Value = -770.5
$ws.Range("H74").Value = 5333.3335
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("H77").Value = 5333.3335
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("H103").Value = 19231276
$ws.Range("J103").Value = 21739660
$ws.Range("L103").Value = 65218980
$ws.Range("N103").Value = -65220152
$ws.Range("H107").Value = 1027
$ws.Range("I107").Value = 930.44446
$ws.Range("J107").Value = 1896
$ws.Range("K107").Value = 930.44446
$ws.Range("L107").Value = 1896
$ws.Range("M107").Value = 989.55554
$ws.Range("N107").Value = -5736
$ws.Range("H113").Value = 4799.077
$ws.Range("I113").Value = 4800
$ws.Range("J113").Value = 4799
$ws.Range("K113").Value = 4800
$ws.Range("L113").Value = 4799
$ws.Range("M113").Value = -1546
$ws.Range("N113").Value = -11307
$ws.Range("H132").Value = 2299.8064
$ws.Range("I132").Value = 1320
$ws.Range("J132").Value = 8913.5
$ws.Range("K132").Value = 3960
$ws.Range("L132").Value = 26740.5
$ws.Range("M132").Value = -1430
$ws.Range("N132").Value = -31800.5
$ws.Range("H133").Value = 85367.60000000001
$ws.Range("J133").Value = 85782.25
$ws.Range("L133").Value = 85782.25
$ws.Range("N133").Value = -95902.25
$ws.Range("H138").Value = 2156.8108
$ws.Range("I138").Value = 1323.381
$ws.Range("J138").Value = 3250.6875
$ws.Range("K138").Value = 3970.143
$ws.Range("L138").Value = 9752.0625
$ws.Range("M138").Value = 1169.857
$ws.Range("N138").Value = -20032.0625
$ws.Range("H141").Value = 41680.652
$ws.Range("I141").Value = 52906.055
$ws.Range("K141").Value = 158718.165
$ws.Range("M141").Value = -153538.165
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 77172.42999999999
$ws.Range("I32").Value = 47947.227
$ws.Range("K32").Value = 47947.227
$ws.Range("M32").Value = -47660.227
$ws.Range("H51").Value = 33999.2
$ws.Range("J51").Value = 33999.2
$ws.Range("L51").Value = 33999.2
$ws.Range("N51").Value = -35511.2
$ws.Range("H74").Value = 1307.7587
$ws.Range("I74").Value = 1270.3478
$ws.Range("J74").Value = 1451.1666
$ws.Range("K74").Value = 1270.3478
$ws.Range("L74").Value = 1451.1666
$ws.Range("M74").Value = -396.3478
$ws.Range("N74").Value = -3199.1666
$ws.Range("H77").Value = 1307.7587
$ws.Range("I77").Value = 1270.3478
$ws.Range("J77").Value = 1451.1666
$ws.Range("K77").Value = 6351.739
$ws.Range("L77").Value = 7255.833000000001
$ws.Range("M77").Value = -1983.739
$ws.Range("N77").Value = -15991.833

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8854
$ws.Range("I20").Value = 7755.7144
$ws.Range("K20").Value = 7755.7144
$ws.Range("M20").Value = -7508.7144

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8246.25
$ws.Range("I86").Value = 7911.6665
$ws.Range("J86").Value = 9250
$ws.Range("K86").Value = 7911.6665
$ws.Range("L86").Value = 9250
$ws.Range("M86").Value = -6788.6665
$ws.Range("N86").Value = -11496
$ws.Range("H89").Value = 8246.25
$ws.Range("I89").Value = 7911.6665
$ws.Range("J89").Value = 9250
$ws.Range("K89").Value = 39558.3325
$ws.Range("L89").Value = 46250
$ws.Range("M89").Value = -33942.3325
$ws.Range("N89").Value = -57482
$ws.Range("H94").Value = 1161.0714
$ws.Range("I94").Value = 1213.4445
$ws.Range("K94").Value = 1213.4445
$ws.Range("M94").Value = -762.4445000000001
$ws.Range("H134").Value = 2545.4285
$ws.Range("I134").Value = 2513.1875
$ws.Range("J134").Value = 2648.6
$ws.Range("K134").Value = 7539.5625
$ws.Range("L134").Value = 7945.799999999999
$ws.Range("M134").Value = -5004.5625
$ws.Range("N134").Value = -13015.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15861662
$ws.Range("I4").Value = 4623166
$ws.Range("K4").Value = 13869498
$ws.Range("M4").Value = -13869386
$ws.Range("H88").Value = 10999.6
$ws.Range("J88").Value = 10999.6
$ws.Range("L88").Value = 32998.8
$ws.Range("N88").Value = -33854.8
$ws.Range("H91").Value = 10999.6
$ws.Range("J91").Value = 10999.6
$ws.Range("L91").Value = 32998.8
$ws.Range("N91").Value = -35962.8
$ws.Range("H113").Value = 449.8
$ws.Range("I113").Value = 466
$ws.Range("J113").Value = 443.9091
$ws.Range("K113").Value = 1398
$ws.Range("L113").Value = 1331.7273
$ws.Range("M113").Value = 772
$ws.Range("N113").Value = -5671.7273
$ws.Range("H128").Value = 347901.5
$ws.Range("I128").Value = 347901.5
$ws.Range("K128").Value = 1043704.5
$ws.Range("M128").Value = -1038724.5
$ws.Range("H140").Value = 3052.9473
$ws.Range("I140").Value = 2039.3529
$ws.Range("J140").Value = 11668.5
$ws.Range("K140").Value = 6118.0587
$ws.Range("L140").Value = 35005.5
$ws.Range("M140").Value = -938.0587000000005
$ws.Range("N140").Value = -45365.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 34966.5
$ws.Range("J52").Value = 34966.5
$ws.Range("L52").Value = 34966.5
$ws.Range("N52").Value = -35484.5
$ws.Range("H70").Value = 6325
$ws.Range("I70").Value = 5276.5
$ws.Range("K70").Value = 5276.5
$ws.Range("M70").Value = -5006.5
$ws.Range("H73").Value = 6325
$ws.Range("I73").Value = 5276.5
$ws.Range("K73").Value = 5276.5
$ws.Range("M73").Value = -4340.5
$ws.Range("H132").Value = 2416.8965
$ws.Range("J132").Value = 4904.6665
$ws.Range("L132").Value = 14713.9995
$ws.Range("N132").Value = -19773.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 130000
$ws.Range("I7").Value = 130000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 130000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -129888
$ws.Range("H16").Value = 860.95
$ws.Range("I16").Value = 914.1875
$ws.Range("K16").Value = 914.1875
$ws.Range("M16").Value = -744.1875
$ws.Range("H32").Value = 1351.9166
$ws.Range("I32").Value = 929.4545000000001
$ws.Range("J32").Value = 5999
$ws.Range("K32").Value = 929.4545000000001
$ws.Range("L32").Value = 5999
$ws.Range("M32").Value = -612.4545000000001
$ws.Range("N32").Value = -6633
$ws.Range("H68").Value = 2450.85
$ws.Range("I68").Value = 2053.4614
$ws.Range("K68").Value = 2053.4614
$ws.Range("M68").Value = -1304.4614
$ws.Range("H71").Value = 2450.85
$ws.Range("I71").Value = 2053.4614
$ws.Range("K71").Value = 10267.307
$ws.Range("M71").Value = -6523.307000000001
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H126").Value = 130000
$ws.Range("I126").Value = 130000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 390000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -387530
$ws.Range("N7").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("N126").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 189999
$ws.Range("J15").Value = 189999
$ws.Range("L15").Value = 189999
$ws.Range("N15").Value = -190575
$ws.Range("H74").Value = 27155.2
$ws.Range("I74").Value = 25500
$ws.Range("K74").Value = 25500
$ws.Range("M74").Value = -24564
$ws.Range("H77").Value = 27155.2
$ws.Range("I77").Value = 25500
$ws.Range("K77").Value = 76500
$ws.Range("M77").Value = -71820
$ws.Range("H122").Value = 5144.6113
$ws.Range("I122").Value = 4850.1875
$ws.Range("K122").Value = 14550.5625
$ws.Range("M122").Value = -12100.5625
$ws.Range("H126").Value = 13512.4
$ws.Range("I126").Value = 15609.3125
$ws.Range("J126").Value = 5124.75
$ws.Range("K126").Value = 46827.9375
$ws.Range("L126").Value = 15374.25
$ws.Range("M126").Value = -44357.9375
$ws.Range("N126").Value = -20314.25
$ws.Range("H132").Value = 37497.43
$ws.Range("I132").Value = 36497.42
$ws.Range("K132").Value = 109492.26
$ws.Range("M132").Value = -106962.26
$ws.Range("H136").Value = 2484.0952
$ws.Range("I136").Value = 2537.6487
$ws.Range("J136").Value = 2087.8
$ws.Range("K136").Value = 7612.946100000001
$ws.Range("L136").Value = 6263.400000000001
$ws.Range("M136").Value = -5062.946100000001
$ws.Range("N136").Value = -11363.4

Write-Host "Applied all cell updates."